$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Marzo de 2020 a las 15:16"

# --- Refresh per-country COVID figures (columns: B=Casos totales, C=Nuevos casos,
#     D=Casos activos, E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes) ---

# Alemania (row 8)
$ws.Cells.Item(8,2).Value  = 17653
$ws.Cells.Item(8,3).Value  = 2333
$ws.Cells.Item(8,5).Value  = 17425
$ws.Cells.Item(8,7).Value  = 4
$ws.Cells.Item(8,8).Value  = 48

# Estados Unidos (row 9)
$ws.Cells.Item(9,2).Value  = 14373
$ws.Cells.Item(9,3).Value  = 584
$ws.Cells.Item(9,5).Value  = 14030

# Suiza (row 12)
$ws.Cells.Item(12,2).Value = 4905
$ws.Cells.Item(12,3).Value = 683
$ws.Cells.Item(12,7).Value = 7
$ws.Cells.Item(12,8).Value = 50

# Austria (row 15)
$ws.Cells.Item(15,2).Value = 2388
$ws.Cells.Item(15,3).Value = 209
$ws.Cells.Item(15,5).Value = 2373

# Noruega (row 17)
$ws.Cells.Item(17,2).Value = 1848
$ws.Cells.Item(17,3).Value = 58
$ws.Cells.Item(17,5).Value = 1840

# Luxemburgo (row 30)
$ws.Cells.Item(30,5).Value = 473
$ws.Cells.Item(30,7).Value = 1
$ws.Cells.Item(30,8).Value = 5

# Singapur (row 37)
$ws.Cells.Item(37,4).Value = 131
$ws.Cells.Item(37,5).Value = 254

# Rusia overtakes Peru / Filipinas / India in the ranking -> rows 50-53 shift.
# Row 50 now holds Rusia's fresh numbers; rows 51-53 keep the untouched
# Peru / Filipinas / India figures, just shifted one row down.
$ws.Cells.Item(50,1).Value = "Rusia"
$ws.Cells.Item(50,2).Value = 253
$ws.Cells.Item(50,3).Value = 54
$ws.Cells.Item(50,4).Value = 12
$ws.Cells.Item(50,5).Value = 240
$ws.Cells.Item(50,6).Value = 0
$ws.Cells.Item(50,7).Value = 0
$ws.Cells.Item(50,8).Value = 1

$ws.Cells.Item(51,1).Value = "Peru"
$ws.Cells.Item(51,2).Value = 234
$ws.Cells.Item(51,3).Value = 0
$ws.Cells.Item(51,4).Value = 1
$ws.Cells.Item(51,5).Value = 230
$ws.Cells.Item(51,6).Value = 7
$ws.Cells.Item(51,7).Value = 2
$ws.Cells.Item(51,8).Value = 3

$ws.Cells.Item(52,1).Value = "Filipinas"
$ws.Cells.Item(52,2).Value = 230
$ws.Cells.Item(52,3).Value = 13
$ws.Cells.Item(52,4).Value = 8
$ws.Cells.Item(52,5).Value = 204
$ws.Cells.Item(52,6).Value = 1
$ws.Cells.Item(52,7).Value = 1
$ws.Cells.Item(52,8).Value = 18

$ws.Cells.Item(53,1).Value = "India"
$ws.Cells.Item(53,2).Value = 223
$ws.Cells.Item(53,3).Value = 29
$ws.Cells.Item(53,4).Value = 23
$ws.Cells.Item(53,5).Value = 195
$ws.Cells.Item(53,6).Value = 0
$ws.Cells.Item(53,7).Value = 1
$ws.Cells.Item(53,8).Value = 5

# Eslovaquia (row 66)
$ws.Cells.Item(66,2).Value = 123
$ws.Cells.Item(66,5).Value = 123
